# Apply updated cryptocurrency market data to the worksheet.
# Numeric-looking text values in column D are prefixed with a leading
# apostrophe so Excel stores them as text (matching the original inlineStr
# cell type) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.943.59"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "3.251.17"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'545.26"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "'147.83"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'0.116"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "'0.433"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "3.808.00"
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "'26.40"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "60.919.51"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "3.255.08"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "'6.32"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "'13.45"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "'8.45"
$ws.Range("E20").Value = "  +2.76%  "
$ws.Range("D21").Value = "'379.21"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'0.533"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "'70.13"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'8.66"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "0.0₃0913"
$ws.Range("E28").Value = "  +3.84%  "
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").Value = "'22.64"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "'6.21"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("D35").Value = "'159.44"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  +6.29%  "
$ws.Range("D37").Value = "'26.38"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("D38").Value = "2.802.57"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("D39").Value = "'0.0722"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").Value = "'0.0313"
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'40.13"
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").Value = "'0.732"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "3.292.36"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("E46").Value = "  +3.19%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'21.56"
$ws.Range("E47").Value = "  +5.03%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'0.806"
$ws.Range("E50").Value = "  +5.72%  "
$ws.Range("D51").Value = "'278.91"
$ws.Range("E51").Value = "  +7.78%  "
